# Updated cryptos list on Wed Jan 31 03:50:31 UTC 2024 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) values for the
# crypto tracker sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that LOOKS numeric (e.g. "306.90") while keeping it
# stored as literal text, matching the sheet's existing inline-string cells
# (avoids Excel's automatic text->number coercion, which would drop
# formatting such as trailing zeros).
function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range("D2").Value = '42.968.88'
$ws.Range("E2").Value = '  -1.25%  '
$ws.Range("D3").Value = '2.342.90'
$ws.Range("E3").Value = '  +1.31%  '
$ws.Range("E4").Value = '  +0.01%  '
Set-TextValue "D5" '306.90'
Set-TextValue "D6" '101.14'
$ws.Range("E6").Value = '  -1.23%  '
Set-TextValue "D7" '0.511'
$ws.Range("E7").Value = '  -4.85%  '
Set-TextValue "D9" '0.513'
$ws.Range("E9").Value = '  -3.56%  '
Set-TextValue "D10" '34.99'
$ws.Range("E10").Value = '  -2.15%  '
$ws.Range("E11").Value = '  +0.12%  '
$ws.Range("E12").Value = '  -2.00%  '
$ws.Range("E13").Value = '  -0.30%  '
$ws.Range("E14").Value = '  -3.01%  '
Set-TextValue "D15" '15.90'
$ws.Range("E15").Value = '  +6.22%  '
$ws.Range("D16").Value = '2.336.23'
$ws.Range("E16").Value = '  +1.10%  '
Set-TextValue "D17" '0.812'
$ws.Range("E17").Value = '  +0.36%  '
$ws.Range("D18").Value = '42.876.34'
$ws.Range("E18").Value = '  -1.24%  '
Set-TextValue "D19" '6.25'
$ws.Range("E19").Value = '  +0.96%  '
$ws.Range("E20").Value = '  -1.71%  '
Set-TextValue "D21" '11.70'
$ws.Range("E21").Value = '  -5.98%  '
Set-TextValue "D22" '67.87'
$ws.Range("E22").Value = '  -0.54%  '
Set-TextValue "D23" '237.14'
$ws.Range("E23").Value = '  -1.93%  '
Set-TextValue "D24" '2.03'
$ws.Range("E24").Value = '  -0.22%  '
$ws.Range("E25").Value = '  -2.47%  '
Set-TextValue "D26" '1.00'
$ws.Range("E26").Value = '  -0.21%  '
Set-TextValue "D27" '25.59'
$ws.Range("E27").Value = '  +3.07%  '
$ws.Range("E28").Value = '  +1.09%  '
Set-TextValue "D29" '35.16'
$ws.Range("E29").Value = '  -4.30%  '
Set-TextValue "D30" '9.37'
$ws.Range("E30").Value = '  -2.69%  '
Set-TextValue "D31" '160.11'
$ws.Range("E31").Value = '  -4.63%  '
Set-TextValue "D32" '0.999'
$ws.Range("E32").Value = '  -0.06%  '
Set-TextValue "D33" '5.14'
$ws.Range("E33").Value = '  -2.97%  '
$ws.Range("E34").Value = '  +8.75%  '
Set-TextValue "D35" '2.49'
$ws.Range("E35").Value = '  -0.60%  '
Set-TextValue "D36" '17.47'
$ws.Range("E36").Value = '  -0.51%  '
Set-TextValue "D37" '0.0729'
Set-TextValue "D38" '2.98'
$ws.Range("E38").Value = '  -3.97%  '
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("E40").Value = '  -3.21%  '
$ws.Range("E41").Value = '  -2.40%  '
$ws.Range("D42").Value = '2.021.06'
$ws.Range("E42").Value = '  +2.54%  '
Set-TextValue "D44" '18.72'
$ws.Range("E44").Value = '  -3.56%  '
Set-TextValue "D45" '10.28'
$ws.Range("E45").Value = '  +3.17%  '
$ws.Range("E46").Value = '  -0.94%  '
Set-TextValue "D47" '56.42'
$ws.Range("E47").Value = '  +1.50%  '
$ws.Range("E48").Value = '  -1.06%  '
$ws.Range("D49").Value = '2.567.40'
$ws.Range("E49").Value = '  +1.11%  '
$ws.Range("E50").Value = '  +2.04%  '
Set-TextValue "D51" '1.52'
$ws.Range("E51").Value = '  -3.41%  '
